$wb = $excel.ActiveWorkbook
$wsOutput = $wb.Worksheets.Item("Output")

# Update the "verify" label on the Output sheet to the more specific
# "verifyclient" wording, and drop the wrap-text formatting that was
# previously applied to that header row.
$wsOutput.Range("A1").Value = "verifyclient"
$wsOutput.Range("A1:B1").WrapText = $false

# Make the Output sheet the active / selected tab, with B1 selected,
# mirroring the periodic & upfront scenario additions.
$wsOutput.Activate() | Out-Null
$wsOutput.Range("B1").Select() | Out-Null
